$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A24").Value = 1099
$ws.Range("B24").Value = "ERROR_TEST"
$ws.Range("C24").Value = "general"
$ws.Range("D24").Value = "测试用错误码"

$ws.Range("E23").Select()
